$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Copy cell formatting from an existing "header + Tiempo1/2/3 + SUMA"
#    block (B41:I45, which carries the s=9/10/6/7/8/3/4/5/11/1 styles we
#    need) down onto the two new blocks being appended to the sheet.
# ------------------------------------------------------------------
$fmtSrc = $ws.Range("B41:I45")
$fmtSrc.Copy()
$ws.Range("B48:I52").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B54:I58").PasteSpecial(-4122)   # xlPasteFormats

# Matching format for the K (media/average) column.
$kSrc = $ws.Range("K45")
$kSrc.Copy()
$ws.Range("K52").PasteSpecial(-4122)
$ws.Range("K58").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Block 1 — "Medición tiempos inserción post índices" (rows 48-52)
# ------------------------------------------------------------------
$ws.Range("B48").Value = "Medición tiempos inserción post índices"
$ws.Range("C48").Value = 4
$ws.Range("D48").Value = "6A"
$ws.Range("E48").Value = "6B"
$ws.Range("F48").Value = "5.1"

$ws.Range("B49").Value = "Tiempo 1"
$ws.Range("C49").Value = 0.013
$ws.Range("D49").Value = 0.017
$ws.Range("E49").Value = 0.021
$ws.Range("F49").Value = 0.07

$ws.Range("B50").Value = "Tiempo 2"
$ws.Range("C50").Value = 0.0052
$ws.Range("D50").Value = 0.0093
$ws.Range("E50").Value = 0.02
$ws.Range("F50").Value = 0.05

$ws.Range("B51").Value = "Tiempo 3"
$ws.Range("C51").Value = 0.0053
$ws.Range("D51").Value = 0.0087
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0.93

$ws.Range("B52").Value = "SUMA"
$ws.Range("C52").Formula = "=SUM(C49:C51)"
$ws.Range("D52:H52").Formula = "=SUM(D49:D51)"
$ws.Range("G52:H52").ClearContents()
$ws.Range("K52").Formula = "=AVERAGE(C52:I52)"

# ------------------------------------------------------------------
# 3. Block 2 — "Tiempos de medición insercción pre índices" (rows 54-58)
# ------------------------------------------------------------------
$ws.Range("B54").Value = "Tiempos de medición insercción pre índices"
$ws.Range("C54").Value = 4
$ws.Range("D54").Value = 5
$ws.Range("E54").Value = "6A"
$ws.Range("F54").Value = "6B"
$ws.Range("G54").Value = "5.2"
$ws.Range("H54").Value = "5.3"

$ws.Range("B55").Value = "Tiempo 1"
$ws.Range("C55").Value = 0.01
$ws.Range("D55").Value = 0.64
$ws.Range("E55").Value = 0.268
$ws.Range("F55").Value = 0.346

$ws.Range("B56").Value = "Tiempo 2"
$ws.Range("C56").Value = 0.01
$ws.Range("D56").Value = 0.06
$ws.Range("E56").Value = 0.011
$ws.Range("F56").Value = 0.317

$ws.Range("B57").Value = "Tiempo 3"
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 0.42
$ws.Range("E57").Value = 0.168

$ws.Range("B58").Value = "Media"
$ws.Range("C58").Formula = "=SUM(C55:C57)"
$ws.Range("D58:I58").Formula = "=SUM(D55:D57)"
$ws.Range("K58").Formula = "=AVERAGE(C58:I58)"

# ------------------------------------------------------------------
# 4. Column widths for the newly-populated C and K columns.
# ------------------------------------------------------------------
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(11).EntireColumn.AutoFit() | Out-Null

# ------------------------------------------------------------------
# 5. Selection mirrors where editing finished.
# ------------------------------------------------------------------
$ws.Range("F57").Select() | Out-Null
